$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a number by Excel
# (losing formatting like trailing zeros) must be forced to stay text by
# setting an explicit text NumberFormat before assigning the value.
$textForceCells = @("D5", "D6", "D9", "D10", "D16", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D30", "D32", "D33", "D34", "D35", "D41", "D42", "D44", "D46", "D48", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row (price + volume refresh; rows 24/25
# additionally swap BitcoinCash <-> PancakeSwap).
$ws.Range("D2").Value = "51.107.37"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.964.19"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "379.53"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "102.27"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "36.56"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "3.432.86"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("E14").Value = "  +6.34%  "
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("D16").Value = "12.00"
$ws.Range("E16").Value = "  +67.28%  "
$ws.Range("D17").Value = "2.974.85"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").Value = "51.222.38"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "3.13"
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "70.12"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "3.28"
$ws.Range("E24").Value = "  +13.32%  "
$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D25").Value = "267.59"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("D27").Value = "7.18"
$ws.Range("E27").Value = "  -8.99%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "25.90"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").Value = "10.38"
$ws.Range("E32").Value = "  +5.21%  "
$ws.Range("D33").Value = "34.50"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").Value = "51.20"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "2.02"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +8.71%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").Value = "16.56"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("D42").Value = "125.40"
$ws.Range("E42").Value = "  +4.17%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").Value = "21.74"
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("E45").Value = "  +9.43%  "
$ws.Range("D46").Value = "2.40"
$ws.Range("E46").Value = "  +4.07%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "0.271"
$ws.Range("E48").Value = "  -7.37%  "
$ws.Range("D49").Value = "2.053.00"
$ws.Range("E49").Value = "  +3.83%  "
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("D51").Value = "5.43"
$ws.Range("E51").Value = "  +7.79%  "
